$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "MCH220-1"
$ws.Range("C2").Value = "MISCELLANEOUS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24E | GRAP COUNT NUMER: NONE"

# --- Row 3 ---
$ws.Range("A3").Value = "MCH220-2"
$ws.Range("C3").Value = "MISCELLANEOUS"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24E | GRAP COUNT NUMER: NONE"

# --- Row 4 ---
$ws.Range("A4").Value = "MCH220-3"
$ws.Range("C4").Value = "MISCELLANEOUS, VARIOUS FILES ON ACTIVITEES"
$ws.Range("E4").Value = "Series"
$ws.Range("F4").Value = "1 Box"
$ws.Range("G4").Value = "LOCATION: 24E | GRAP COUNT NUMER: NONE"

# --- Formatting ---
# (D and H are intentionally left without a .Value write: applying the font
#  below still materialises the cell, but empty, exactly like the source.)
# Most cells in the new rows use a plain Calibri 10 (theme) font.
$bodyRange = $ws.Range("A2:A4,C2:E4,G2:H4")
foreach ($area in $bodyRange.Areas) {
    $area.Font.Name = "Calibri"
    $area.Font.Size = 10
    $area.Font.ThemeColor = 1
}

# The extentAndMedium (F) column gets its own (visually identical) style.
$extentRange = $ws.Range("F2:F4")
$extentRange.Font.Name = "Calibri"
$extentRange.Font.Size = 10
$extentRange.Font.ThemeColor = 1

# --- View state ---
# Keep the header row frozen and move the selection onto the new data.
$win = $excel.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("A2:L4").Select()
